$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 284; $r -le 330; $r++) {
    $ws.Cells.Item($r, 7).Value = 46.8799663
    $ws.Cells.Item($r, 8).Value = -121.7269094
}
